$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.429.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.643.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.34%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.642.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.63%  "

$ws.Range("E10").Value = "  +6.63%  "

$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("E13").Value = "  +1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.24%  "

$ws.Range("E15").Value = "  +3.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.122.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.311.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.645.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "367.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("E23").Value = "  +2.56%  "

$ws.Range("E24").Value = "  +2.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("E28").Value = "  +6.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.770.65"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "574.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "

$ws.Range("E32").Value = "  +4.65%  "

$ws.Range("E33").Value = "  +4.68%  "

$ws.Range("E34").Value = "  +2.61%  "

$ws.Range("E35").Value = "  +3.21%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("E37").Value = "  +3.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.61%  "

$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.367"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.21%  "

$ws.Range("E44").Value = "  +3.84%  "

$ws.Range("E45").Value = "  +13.99%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "157.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "

$ws.Range("E49").Value = "  +0.94%  "

$ws.Range("E50").Value = "  +1.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.39%  "
